# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Thu Aug 10 13:38:07 UTC 2023 with GitHub Actions".
#
# Columns D ("Price") and E ("Volume(1h)") hold plain text values, even
# though many look numeric (e.g. "0.9985", "5.029", "1.000"). Assigning such
# a string straight to Range.Value lets Excel auto-convert it to a real
# number, changing the cell type and silently dropping formatting such as
# trailing zeros ("1.000" -> 1, "24.70" -> 24.7). To keep those cells as text
# we prefix the assignment with a leading single quote (Excel's classic
# "force text" input trick), then reset the cell Style back to "Normal" so
# no left-over quote-prefix/text number-format remains on the cell - matching
# the original (unstyled) text cells exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.569.48'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '1.856.33'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("D4").Value = "'" + '0.9985'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = "'" + '242.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.07%  '
$ws.Range("D6").Value = "'" + '0.6342'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.03%  '
$ws.Range("D7").Value = "'" + '0.9994'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = "'" + '0.07610'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").Value = "'" + '0.2994'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").Value = "'" + '24.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("D11").Value = "'" + '0.07742'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").Value = '1.861.88'
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").Value = "'" + '0.6944'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("D14").Value = "'" + '5.029'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").Value = "'" + '83.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = "'" + '0.00001000'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.26%  '
$ws.Range("D17").Value = '2.115.55'
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").Value = "'" + '6.276'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").Value = '29.582.27'
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("D20").Value = "'" + '234.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.13%  '
$ws.Range("D21").Value = "'" + '12.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("D22").Value = "'" + '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = "'" + '7.673'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.25%  '
$ws.Range("D24").Value = "'" + '0.9990'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("D25").Value = "'" + '155.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.17%  '
$ws.Range("D26").Value = "'" + '0.1402'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.58%  '
$ws.Range("D27").Value = "'" + '8.480'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("D28").Value = "'" + '17.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").Value = "'" + '0.05830'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.19%  '
$ws.Range("D31").Value = "'" + '1.262'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.16%  '
$ws.Range("D32").Value = "'" + '4.137'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.99%  '
$ws.Range("D33").Value = "'" + '4.033'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.16%  '
$ws.Range("D34").Value = "'" + '1.903'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("D35").Value = "'" + '1.170'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("D36").Value = "'" + '0.7219'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.47%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.248.32'
$ws.Range("E38").Value = '  +2.50%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = "'" + '2.804'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("D40").Value = "'" + '0.01807'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("D41").Value = "'" + '0.9080'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.00%  '
$ws.Range("D42").Value = "'" + '6.131'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.18%  '
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").Value = '2.019.20'
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").Value = "'" + '67.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("D46").Value = "'" + '101.48'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").Value = "'" + '7.356'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.40%  '
$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").Value = "'" + '0.4059'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'" + '9.202'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("D50").Value = "'" + '1.715'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.31%  '
$ws.Range("E51").Value = '  -3.57%  '

Write-Host "Updated 103 cells in the cryptos list."
